$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work from the bottom of the sheet upward so row numbers for
# not-yet-processed rows stay stable.

# 1) Remove row 24 (U664 / C15, no comments/resolution) - merged into U663 row in new data
$ws.Rows(24).Delete()

# 2) Remove rows 19-22 (U659, U660, U661, U662 / C15) - these animals were dropped
$ws.Range("A19:D22").EntireRow.Delete()

# 3) Row 17 held U583/C14; QC rerun replaced it with two new animals: U402 (C10) and U459 (C11).
#    First insert a new blank row after row 17 to hold the second new record...
$ws.Rows(18).Insert()

# ...then populate row 17 with the first new record (U402) including its resolution notes,
# and row 18 with the second new record (U459).
$ws.Range("A17").Value = "U402"
$ws.Range("B17").Value = "C10"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = "EXCLUDE_LOCOMOTOR2, EXCLUDE_DELAYED_PUNISHMENT"

$ws.Range("A18").Value = "U459"
$ws.Range("B18").Value = "C11"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = ""

# 4) Remove row 2 (U3 / C01) - animal dropped entirely from the QC'd dataset
$ws.Rows(2).Delete()
